$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-19 Friday" "2024-07-20 Saturday"

Replace-Text "750÷3=250, 0" "168÷7=24, 0"
Replace-Text "815÷6=135, 5" "195÷8=24, 3"
Replace-Text "952÷3=317, 1" "140÷3=46, 2"
Replace-Text "611÷6=101, 5" "658÷6=109, 4"
Replace-Text "184÷2=92, 0" "849÷5=169, 4"

Replace-Text "369÷4=92, 1" "803÷9=89, 2"
Replace-Text "814÷5=162, 4" "616÷2=308, 0"
Replace-Text "542÷7=77, 3" "961÷5=192, 1"
Replace-Text "745÷5=149, 0" "554÷9=61, 5"
Replace-Text "260÷4=65, 0" "330÷5=66, 0"

Replace-Text "771÷4=192, 3" "946÷7=135, 1"
Replace-Text "682÷2=341, 0" "592÷6=98, 4"
Replace-Text "140÷5=28, 0" "156÷7=22, 2"
Replace-Text "553÷7=79, 0" "338÷6=56, 2"
Replace-Text "322÷7=46, 0" "894÷9=99, 3"

Replace-Text "111÷9=12, 3" "142÷8=17, 6"
Replace-Text "715÷3=238, 1" "395÷9=43, 8"
Replace-Text "388÷2=194, 0" "552÷9=61, 3"
Replace-Text "267÷9=29, 6" "832÷5=166, 2"
Replace-Text "379÷6=63, 1" "127÷6=21, 1"

Replace-Text "861÷2=430, 1" "644÷4=161, 0"
Replace-Text "991÷6=165, 1" "788÷5=157, 3"
Replace-Text "991÷9=110, 1" "887÷2=443, 1"
Replace-Text "186÷2=93, 0" "648÷3=216, 0"
Replace-Text "387÷3=129, 0" "675÷6=112, 3"
